$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$images = @(
    "Celebs/14_Emma_Watson_0004.jpg",
    "Celebs/14_Emma_Watson_0005.jpg",
    "Celebs/14_Emma_Watson_0008.jpg",
    "Celebs/14_Emma_Watson_0013.jpg",
    "Celebs/14_Emma_Watson_0015.jpg",
    "Celebs/15_Margot_Robbie_0005.jpg",
    "Celebs/15_Margot_Robbie_0012.jpg",
    "Celebs/16_Emma_Stone_0008.jpg",
    "Celebs/16_Emma_Stone_0010.jpg",
    "Celebs/16_Emma_Stone_0012.jpg",
    "Celebs/16_Emma_Stone_0014.jpg",
    "Celebs/16_Margot_Robbie_0009.jpg",
    "Celebs/17_Daniel_Radcliffe_0007.jpg",
    "Celebs/17_Emma_Watson_0006.jpg",
    "Celebs/17_Margot_Robbie_0006.jpg",
    "Celebs/17_Margot_Robbie_0012.jpg",
    "Celebs/17_Margot_Robbie_0014.jpg",
    "Celebs/18_Daniel_Radcliffe_0002.jpg",
    "Celebs/18_Daniel_Radcliffe_0007.jpg",
    "Celebs/18_Daniel_Radcliffe_0008.jpg",
    "Celebs/18_Emma_Stone_0016.jpg",
    "Celebs/18_Emma_Watson_0004.jpg",
    "Celebs/19_Emma_Stone_0006.jpg",
    "Celebs/19_Margot_Robbie_0001.jpg",
    "Celebs/19_Margot_Robbie_0004.jpg",
    "Celebs/20_Emma_Stone_0005.jpg",
    "Celebs/20_Emma_Watson_0007.jpg",
    "Celebs/20_Margot_Robbie_0004.jpg",
    "Celebs/20_Margot_Robbie_0006.jpg",
    "Celebs/20_Scarlett_Johansson_0011.jpg",
    "Celebs/21_Daniel_Radcliffe_0002.jpg",
    "Celebs/21_Daniel_Radcliffe_0008.jpg",
    "Celebs/21_Daniel_Radcliffe_0009.jpg",
    "Celebs/21_Daniel_Radcliffe_0012.jpg",
    "Celebs/21_Emma_Stone_0002.jpg",
    "Celebs/21_Emma_Stone_0003.jpg",
    "Celebs/21_Emma_Watson_0004.jpg",
    "Celebs/21_Margot_Robbie_0003.jpg",
    "Celebs/21_Scarlett_Johansson_0007.jpg",
    "Celebs/22_Emma_Stone_0013.jpg",
    "Celebs/22_Emma_Watson_0001.jpg",
    "Celebs/22_Emma_Watson_0003.jpg",
    "Celebs/22_Emma_Watson_0010.jpg",
    "Celebs/22_Margot_Robbie_0010.jpg",
    "Celebs/22_Scarlett_Johansson_0005.jpg",
    "Celebs/23_Daniel_Radcliffe_0002.jpg",
    "Celebs/23_Daniel_Radcliffe_0010.jpg",
    "Celebs/23_Daniel_Radcliffe_0019.jpg",
    "Celebs/23_Emma_Stone_0014.jpg",
    "Celebs/24_Daniel_Radcliffe_0017.jpg",
    "Celebs/24_Emma_Stone_0015.jpg",
    "Celebs/24_Scarlett_Johansson_0001.jpg",
    "Celebs/24_Scarlett_Johansson_0011.jpg",
    "Celebs/25_Scarlett_Johansson_0010.jpg",
    "Celebs/26_Scarlett_Johansson_0003.jpg",
    "Celebs/26_Scarlett_Johansson_0008.jpg",
    "Celebs/27_Scarlett_Johansson_0012.jpg",
    "Celebs/28_Scarlett_Johansson_0005.jpg",
    "Celebs/29_Scarlett_Johansson_0010.jpg",
    "Celebs/29_Scarlett_Johansson_0016.jpg",
    "Celebs/31_Neil_Patrick_Harris_0008.jpg",
    "Celebs/32_Jim_Parsons_0010.jpg",
    "Celebs/32_Jim_Parsons_0016.jpg",
    "Celebs/32_Jim_Parsons_0018.jpg",
    "Celebs/33_Jim_Parsons_0011.jpg",
    "Celebs/33_Neil_Patrick_Harris_0011.jpg",
    "Celebs/34_Neil_Patrick_Harris_0004.jpg",
    "Celebs/34_Neil_Patrick_Harris_0007.jpg",
    "Celebs/35_Neil_Patrick_Harris_0009.jpg",
    "Celebs/36_Jim_Parsons_0011.jpg",
    "Celebs/36_Neil_Patrick_Harris_0007.jpg",
    "Celebs/36_Neil_Patrick_Harris_0010.jpg",
    "Celebs/37_Neil_Patrick_Harris_0010.jpg",
    "Celebs/37_Neil_Patrick_Harris_0014.jpg",
    "Celebs/37_Will_Smith_0015.jpg",
    "Celebs/38_Jim_Parsons_0003.jpg",
    "Celebs/38_Neil_Patrick_Harris_0006.jpg",
    "Celebs/38_Will_Smith_0013.jpg",
    "Celebs/39_Jim_Parsons_0001.jpg",
    "Celebs/39_Jim_Parsons_0004.jpg",
    "Celebs/39_Jim_Parsons_0008.jpg",
    "Celebs/39_Jim_Parsons_0010.jpg",
    "Celebs/39_Neil_Patrick_Harris_0001.jpg",
    "Celebs/39_Neil_Patrick_Harris_0010.jpg",
    "Celebs/40_Jim_Parsons_0004.jpg",
    "Celebs/40_Jim_Parsons_0015.jpg",
    "Celebs/40_Will_Smith_0012.jpg",
    "Celebs/42_Will_Smith_0012.jpg",
    "Celebs/43_Will_Smith_0006.jpg",
    "Celebs/43_Will_Smith_0007.jpg",
    "Celebs/43_Will_Smith_0010.jpg",
    "Celebs/44_Will_Smith_0006.jpg",
    "Celebs/44_Will_Smith_0007.jpg",
    "Celebs/44_Will_Smith_0015.jpg",
    "Celebs/45_Johnny_Depp_0003.jpg",
    "Celebs/45_Will_Smith_0010.jpg",
    "Celebs/45_Will_Smith_0014.jpg",
    "Celebs/47_Johnny_Depp_0003.jpg",
    "Celebs/47_Johnny_Depp_0004.jpg",
    "Celebs/47_Johnny_Depp_0009.jpg",
    "Celebs/48_Johnny_Depp_0001.jpg",
    "Celebs/48_Johnny_Depp_0003.jpg",
    "Celebs/48_Johnny_Depp_0005.jpg",
    "Celebs/49_Johnny_Depp_0013.jpg",
    "Celebs/50_Jackie_Chan_0012.jpg",
    "Celebs/50_Jackie_Chan_0014.jpg",
    "Celebs/50_Johnny_Depp_0004.jpg",
    "Celebs/50_Johnny_Depp_0008.jpg",
    "Celebs/50_Johnny_Depp_0010.jpg",
    "Celebs/50_Johnny_Depp_0018.jpg",
    "Celebs/51_Jackie_Chan_0006.jpg",
    "Celebs/52_Jackie_Chan_0002.jpg",
    "Celebs/54_Jackie_Chan_0001.jpg",
    "Celebs/54_Jackie_Chan_0014.jpg",
    "Celebs/55_Jackie_Chan_0003.jpg",
    "Celebs/55_Jackie_Chan_0008.jpg",
    "Celebs/55_Jackie_Chan_0015.jpg",
    "Celebs/57_Jackie_Chan_0008.jpg",
    "Celebs/57_Jackie_Chan_0010.jpg",
    "Celebs/58_Jackie_Chan_0005.jpg"
)

$labels = @(
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "target",
    "target",
    "target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "target",
    "target",
    "target",
    "target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "target",
    "target",
    "target",
    "non_target",
    "target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target",
    "non_target"
)

# Header row
$ws.Cells.Item(1, 1).Value = "image"
$ws.Cells.Item(1, 2).Value = "stimulus"

# Copy the header style (bold, border, centered) from A1 to B1
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $images.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $images[$i]
    $ws.Cells.Item($row, 2).Value = $labels[$i]
}

Write-Host "Updated sheet with" $images.Length "data rows. Used range:" $ws.UsedRange.Address()
